$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap/rotate corrected match data in existing rows (home/away team mixups) ---
# Row 16 <= corrected data (was row 17)
$ws.Range("F16").Value = "Mornar Bar"
$ws.Range("G16").Value = 0
$ws.Range("H16").Value = "Arsenal Tivat"
$ws.Range("I16").Value = 1
$ws.Range("J16").Value = 2.5
$ws.Range("K16").Value = "12/08/2023 17:12"
$ws.Range("L16").Value = 2.64
$ws.Range("M16").Value = "12/08/2023 20:14"
$ws.Range("N16").Value = 2.87
$ws.Range("O16").Value = "12/08/2023 17:12"
$ws.Range("P16").Value = 2.89
$ws.Range("Q16").Value = "12/08/2023 18:34"
$ws.Range("R16").Value = 3.03
$ws.Range("S16").Value = "12/08/2023 17:12"
$ws.Range("T16").Value = 2.89
$ws.Range("U16").Value = "12/08/2023 20:14"
$ws.Range("V16").Value = "https://www.betexplorer.com/football/montenegro/prva-crnogorska-liga/mornar-bar-arsenal-tivat/v7NzFNk5/"

# Row 17 <= corrected data (was row 16)
$ws.Range("F17").Value = "Petrovac"
$ws.Range("G17").Value = 2
$ws.Range("H17").Value = "Rudar"
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 1.84
$ws.Range("K17").Value = "11/08/2023 17:42"
$ws.Range("L17").Value = 1.67
$ws.Range("M17").Value = "12/08/2023 20:25"
$ws.Range("N17").Value = 3.27
$ws.Range("O17").Value = "11/08/2023 17:42"
$ws.Range("P17").Value = 3.68
$ws.Range("Q17").Value = "12/08/2023 20:25"
$ws.Range("R17").Value = 3.91
$ws.Range("S17").Value = "11/08/2023 17:42"
$ws.Range("T17").Value = 5.01
$ws.Range("U17").Value = "12/08/2023 20:25"
$ws.Range("V17").Value = "https://www.betexplorer.com/football/montenegro/prva-crnogorska-liga/petrovac-rudar/UeGmCLZN/"

# Row 18 <= corrected data (was row 20)
$ws.Range("F18").Value = "Buducnost"
$ws.Range("G18").Value = 1
$ws.Range("H18").Value = "Jezero"
$ws.Range("I18").Value = 1
$ws.Range("J18").Value = 1.48
$ws.Range("K18").Value = "12/08/2023 09:12"
$ws.Range("L18").Value = 1.47
$ws.Range("M18").Value = "13/08/2023 19:59"
$ws.Range("N18").Value = 3.83
$ws.Range("O18").Value = "12/08/2023 09:12"
$ws.Range("P18").Value = 3.91
$ws.Range("Q18").Value = "13/08/2023 19:59"
$ws.Range("R18").Value = 5.58
$ws.Range("S18").Value = "12/08/2023 09:12"
$ws.Range("T18").Value = 7.31
$ws.Range("U18").Value = "13/08/2023 19:59"
$ws.Range("V18").Value = "https://www.betexplorer.com/football/montenegro/prva-crnogorska-liga/buducnost-jezero/nXQrD1KH/"

# Row 20 <= corrected data (was row 18)
$ws.Range("F20").Value = "Jedinstvo"
$ws.Range("G20").Value = 0
$ws.Range("H20").Value = "Sutjeska"
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 5.84
$ws.Range("K20").Value = "12/08/2023 17:12"
$ws.Range("L20").Value = 4.16
$ws.Range("M20").Value = "13/08/2023 19:20"
$ws.Range("N20").Value = 3.85
$ws.Range("O20").Value = "12/08/2023 17:12"
$ws.Range("P20").Value = 3.62
$ws.Range("Q20").Value = "13/08/2023 19:20"
$ws.Range("R20").Value = 1.48
$ws.Range("S20").Value = "12/08/2023 17:12"
$ws.Range("T20").Value = 1.8
$ws.Range("U20").Value = "13/08/2023 19:20"
$ws.Range("V20").Value = "https://www.betexplorer.com/football/montenegro/prva-crnogorska-liga/jedinstvo-sutjeska/YTMvEs5B/"

# Row 28 <= corrected data (was row 30)
$ws.Range("F28").Value = "Decic"
$ws.Range("G28").Value = 0
$ws.Range("H28").Value = "Rudar"
$ws.Range("I28").Value = 1
$ws.Range("J28").Value = 1.58
$ws.Range("K28").Value = "25/08/2023 08:13"
$ws.Range("L28").Value = 1.57
$ws.Range("M28").Value = "26/08/2023 19:54"
$ws.Range("N28").Value = 3.5
$ws.Range("O28").Value = "25/08/2023 08:13"
$ws.Range("P28").Value = 3.54
$ws.Range("Q28").Value = "26/08/2023 19:58"
$ws.Range("R28").Value = 5.09
$ws.Range("S28").Value = "25/08/2023 08:13"
$ws.Range("T28").Value = 6.07
$ws.Range("U28").Value = "26/08/2023 19:58"
$ws.Range("V28").Value = "https://www.betexplorer.com/football/montenegro/prva-crnogorska-liga/decic-rudar/6u2b7Wmd/"

# Row 29 <= corrected data (was row 28)
$ws.Range("F29").Value = "Buducnost"
$ws.Range("G29").Value = 2
$ws.Range("H29").Value = "Arsenal Tivat"
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 1.4
$ws.Range("K29").Value = "25/08/2023 08:13"
$ws.Range("L29").Value = 1.48
$ws.Range("M29").Value = "26/08/2023 18:05"
$ws.Range("N29").Value = 4.1
$ws.Range("O29").Value = "25/08/2023 08:13"
$ws.Range("P29").Value = 4.06
$ws.Range("Q29").Value = "26/08/2023 18:05"
$ws.Range("R29").Value = 6.27
$ws.Range("S29").Value = "25/08/2023 08:13"
$ws.Range("T29").Value = 6.67
$ws.Range("U29").Value = "26/08/2023 18:05"
$ws.Range("V29").Value = "https://www.betexplorer.com/football/montenegro/prva-crnogorska-liga/buducnost-arsenal-tivat/OIyNQKlH/"

# Row 30 <= corrected data (was row 29)
$ws.Range("F30").Value = "Mornar Bar"
$ws.Range("G30").Value = 4
$ws.Range("H30").Value = "Jedinstvo"
$ws.Range("I30").Value = 3
$ws.Range("J30").Value = 1.75
$ws.Range("K30").Value = "25/08/2023 08:13"
$ws.Range("L30").Value = 2.14
$ws.Range("M30").Value = "26/08/2023 19:58"
$ws.Range("N30").Value = 3.23
$ws.Range("O30").Value = "25/08/2023 08:13"
$ws.Range("P30").Value = 2.89
$ws.Range("Q30").Value = "26/08/2023 19:25"
$ws.Range("R30").Value = 4.3
$ws.Range("S30").Value = "25/08/2023 08:13"
$ws.Range("T30").Value = 3.28
$ws.Range("U30").Value = "26/08/2023 19:58"
$ws.Range("V30").Value = "https://www.betexplorer.com/football/montenegro/prva-crnogorska-liga/mornar-bar-jedinstvo/6ssIR0ZA/"

# Row 37 <= corrected data (was row 38)
$ws.Range("F37").Value = "Mladost DG"
$ws.Range("G37").Value = 2
$ws.Range("H37").Value = "Arsenal Tivat"
$ws.Range("I37").Value = 1
$ws.Range("J37").Value = 2.08
$ws.Range("K37").Value = "15/09/2023 06:43"
$ws.Range("L37").Value = 3.02
$ws.Range("M37").Value = "16/09/2023 19:29"
$ws.Range("N37").Value = 2.91
$ws.Range("O37").Value = "15/09/2023 06:43"
$ws.Range("P37").Value = 2.95
$ws.Range("Q37").Value = "16/09/2023 18:58"
$ws.Range("R37").Value = 3.43
$ws.Range("S37").Value = "15/09/2023 06:43"
$ws.Range("T37").Value = 2.5
$ws.Range("U37").Value = "16/09/2023 19:29"
$ws.Range("V37").Value = "https://www.betexplorer.com/football/montenegro/prva-crnogorska-liga/mladost-dg-arsenal-tivat/naNXLV2k/"

# Row 38 <= corrected data (was row 37)
$ws.Range("F38").Value = "Petrovac"
$ws.Range("G38").Value = 1
$ws.Range("H38").Value = "Jedinstvo"
$ws.Range("I38").Value = 1
$ws.Range("J38").Value = 1.79
$ws.Range("K38").Value = "15/09/2023 06:43"
$ws.Range("L38").Value = 1.8
$ws.Range("M38").Value = "16/09/2023 19:06"
$ws.Range("N38").Value = 3.25
$ws.Range("O38").Value = "15/09/2023 06:43"
$ws.Range("P38").Value = 3.41
$ws.Range("Q38").Value = "16/09/2023 19:06"
$ws.Range("R38").Value = 4.02
$ws.Range("S38").Value = "15/09/2023 06:43"
$ws.Range("T38").Value = 4.5
$ws.Range("U38").Value = "16/09/2023 19:06"
$ws.Range("V38").Value = "https://www.betexplorer.com/football/montenegro/prva-crnogorska-liga/petrovac-jedinstvo/88OTMBmq/"

# --- Append new rows 72-74 (new matches scraped) ---
# Copy row 71 formatting (bold/border style on col A, date style on col E) down to the new rows
$ws.Range("A71:V71").Copy()
$ws.Range("A72:V74").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 72
$ws.Range("A72").Value = 71
$ws.Range("B72").Value = "montenegro"
$ws.Range("C72").Value = "prva-crnogorska-liga"
$ws.Range("D72").Value = "2023-2024"
$ws.Range("E72").Value = 45235.5625
$ws.Range("F72").Value = "Arsenal Tivat"
$ws.Range("G72").Value = 1
$ws.Range("H72").Value = "Buducnost"
$ws.Range("I72").Value = 1
$ws.Range("J72").Value = 4.24
$ws.Range("K72").Value = "04/11/2023 01:43"
$ws.Range("L72").Value = 5.09
$ws.Range("M72").Value = "05/11/2023 06:36"
$ws.Range("N72").Value = 3.29
$ws.Range("O72").Value = "04/11/2023 01:43"
$ws.Range("P72").Value = 3.75
$ws.Range("Q72").Value = "05/11/2023 06:36"
$ws.Range("R72").Value = 1.74
$ws.Range("S72").Value = "04/11/2023 01:43"
$ws.Range("T72").Value = 1.63
$ws.Range("U72").Value = "05/11/2023 06:36"
$ws.Range("V72").Value = "https://www.betexplorer.com/football/montenegro/prva-crnogorska-liga/arsenal-tivat-buducnost/GQFj3Kin/"

# Row 73
$ws.Range("A73").Value = 72
$ws.Range("B73").Value = "montenegro"
$ws.Range("C73").Value = "prva-crnogorska-liga"
$ws.Range("D73").Value = "2023-2024"
$ws.Range("E73").Value = 45235.5625
$ws.Range("F73").Value = "Rudar"
$ws.Range("G73").Value = 0
$ws.Range("H73").Value = "Decic"
$ws.Range("I73").Value = 1
$ws.Range("J73").Value = 4.05
$ws.Range("K73").Value = "04/11/2023 01:43"
$ws.Range("L73").Value = 3.94
$ws.Range("M73").Value = "05/11/2023 13:24"
$ws.Range("N73").Value = 3.2
$ws.Range("O73").Value = "04/11/2023 01:43"
$ws.Range("P73").Value = 3.06
$ws.Range("Q73").Value = "05/11/2023 13:24"
$ws.Range("R73").Value = 1.8
$ws.Range("S73").Value = "04/11/2023 01:43"
$ws.Range("T73").Value = 2.04
$ws.Range("U73").Value = "05/11/2023 13:11"
$ws.Range("V73").Value = "https://www.betexplorer.com/football/montenegro/prva-crnogorska-liga/rudar-decic/UkVe0Iy5/"

# Row 74
$ws.Range("A74").Value = 73
$ws.Range("B74").Value = "montenegro"
$ws.Range("C74").Value = "prva-crnogorska-liga"
$ws.Range("D74").Value = "2023-2024"
$ws.Range("E74").Value = 45235.58333333334
$ws.Range("F74").Value = "Jedinstvo"
$ws.Range("G74").Value = 1
$ws.Range("H74").Value = "Mornar Bar"
$ws.Range("I74").Value = 1
$ws.Range("J74").Value = 2.58
$ws.Range("K74").Value = "04/11/2023 02:13"
$ws.Range("L74").Value = 2.83
$ws.Range("M74").Value = "05/11/2023 13:40"
$ws.Range("N74").Value = 2.77
$ws.Range("O74").Value = "04/11/2023 02:13"
$ws.Range("P74").Value = 2.9
$ws.Range("Q74").Value = "05/11/2023 13:40"
$ws.Range("R74").Value = 2.74
$ws.Range("S74").Value = "04/11/2023 02:13"
$ws.Range("T74").Value = 2.63
$ws.Range("U74").Value = "05/11/2023 13:40"
$ws.Range("V74").Value = "https://www.betexplorer.com/football/montenegro/prva-crnogorska-liga/jedinstvo-mornar-bar/2NBn40xt/"
